# Generate Report for Archive
#
# The localization status report is regenerated: the file
# "18325dcf-0f05-4d53-bbcd-5278d190a6b1.md" has moved out of the
# "Ready for handoff" bucket into the "In Translation" bucket, so the
# report table (present on all three sheets: Overview, zh-cn, de-de)
# is re-sorted -- "In Translation" rows first (alphabetically by file
# name), then "Ready for handoff" rows, with row 2 now holding the
# 18325dcf record, row 3 the ce90b908 record and row 4 the e2c01170
# record. The 3521f6a1 / .localization-config rows are unaffected.
#
# Per worksheet, only the "identity" columns (the source file name,
# and -- on the language sheets -- the handoff file name/date) travel
# with the row; the Status column is recomputed for the new grouping
# (rows 2-4 are now all "In Translation", row 5 stays
# "Ready for handoff"), and the Handback-DateTime / Handoff-Reason
# columns stay put since they are driven by row position in the
# template. Hyperlink targets (r:id -> external URL) are left exactly
# as-is (Excel keeps the relationship bound to the cell address), only
# the visible hyperlink text is refreshed to match the new cell value.

function Set-HyperlinkDisplay($ws, $addr, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $newText
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Overview" -- columns A (File Name), B (zh-cn), C (de-de)
# Only column A travels with the row identity; B/C (status per
# language) stay fixed to the row position.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$a2 = $ws1.Range("A2").Value()
$a3 = $ws1.Range("A3").Value()
$a4 = $ws1.Range("A4").Value()

# new A2 <- old A4, new A3 <- old A2, new A4 <- old A3
$ws1.Range("A2").Value = $a4
$ws1.Range("A3").Value = $a2
$ws1.Range("A4").Value = $a3

# 18325dcf (now row 2) has moved into the "In Translation" bucket
$ws1.Range("B2:C4").Value = "In Translation"

Set-HyperlinkDisplay $ws1 '$A$2' $a4
Set-HyperlinkDisplay $ws1 '$A$3' $a2
Set-HyperlinkDisplay $ws1 '$A$4' $a3

# ---------------------------------------------------------------
# Sheets 2 & 3: "zh-cn" / "de-de" -- columns A (Source File Name),
# C (Latest Handoff File) and D (Latest Handoff Datetime) travel
# together with the row identity; B/G/H stay fixed to the row
# position.
# ---------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $a2 = $ws.Range("A2").Value()
    $c2 = $ws.Range("C2").Value()
    $d2 = $ws.Range("D2").Value()

    $a3 = $ws.Range("A3").Value()
    $c3 = $ws.Range("C3").Value()
    $d3 = $ws.Range("D3").Value()

    $a4 = $ws.Range("A4").Value()
    $c4 = $ws.Range("C4").Value()
    $d4 = $ws.Range("D4").Value()

    # new row2 <- old row4, new row3 <- old row2, new row4 <- old row3
    $ws.Range("A2").Value = $a4
    $ws.Range("C2").Value = $c4
    $ws.Range("D2").Value = $d4

    $ws.Range("A3").Value = $a2
    $ws.Range("C3").Value = $c2
    $ws.Range("D3").Value = $d2

    $ws.Range("A4").Value = $a3
    $ws.Range("C4").Value = $c3
    $ws.Range("D4").Value = $d3

    Set-HyperlinkDisplay $ws '$A$2' $a4
    Set-HyperlinkDisplay $ws '$C$2' $c4
    Set-HyperlinkDisplay $ws '$A$3' $a2
    Set-HyperlinkDisplay $ws '$C$3' $c2
    Set-HyperlinkDisplay $ws '$A$4' $a3
    Set-HyperlinkDisplay $ws '$C$4' $c3
}
